$d = $word.ActiveDocument

# --- First paragraph: add a paragraph border (space-only, no line) and
#     change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1 = $d.Paragraphs.First
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25

# --- First paragraph text: replace the old bookmark-style id text and drop
#     the trailing space run, leaving a single run with the new id text.
$d.Content.Find.Execute("**ID__AFFARS_mp_5315_3_topic_41__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP_5315_3_APPENDIX_A__ID**", 2)
